$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 515, shifting existing rows 515:578 down to 516:579
$ws.Rows.Item(515).Insert()

# Populate the newly inserted row 515 with the new weekly record
$ws.Cells.Item(515, 1).Value2 = 8
$ws.Cells.Item(515, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(515, 3).Value2 = "Coquimbo"
$ws.Cells.Item(515, 4).Value2 = 45142
$ws.Cells.Item(515, 5).Value2 = 4
$ws.Cells.Item(515, 6).Value2 = 100114013
$ws.Cells.Item(515, 7).Value2 = "Zanahoria"
$ws.Cells.Item(515, 8).Value2 = "Sin especificar"
$ws.Cells.Item(515, 9).Value2 = "Primera"
$ws.Cells.Item(515, 10).Value2 = 440
$ws.Cells.Item(515, 11).Value2 = 6000
$ws.Cells.Item(515, 12).Value2 = 6500
$ws.Cells.Item(515, 13).Value2 = 6250
$ws.Cells.Item(515, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(515, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(515, 16).Value2 = 312
$ws.Cells.Item(515, 17).Value2 = 20
$ws.Cells.Item(515, 18).Value2 = "Hortaliza"
